# Project 2 presentation updated
# Adds 5 new submission-log rows (ridge10..ridge14) with updated feature lists.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: ridge10_submit.csv ---
$ws.Cells.Item(11, 1).Value2 = "ridge10_submit.csv"
$ws.Cells.Item(11, 2).Value2 = "ridge"
$ws.Cells.Item(11, 3).Value2 = 1000
$ws.Cells.Item(11, 4).Value2 = "['n_grnhill', 'n_greens', 'n_blueste', 'n_npkvill', 'n_veenker', 'n_brdale', 'n_blmngtn', 'n_meadowv', 'n_clearcr', 'n_swisu', 'n_stonebr', 'n_timber', 'n_noridge', 'n_idotrr', 'n_crawfor', 'n_brkside', 'n_mitchel', 'n_sawyerw', 'n_nwames', 'n_sawyer', 'n_gilbert', 'n_nridght', 'n_somerst', 'n_edwards', 'n_oldtown', 'n_collgcr', 'n_names', 'type_twn', 'type_sf', 'type_twn_end', 'gar_attached', 'gar_detached', 'gar_builtin', 'gar_basement', 'gar_2types', 'gar_carport', 'quality', 'gr_living_sqft', 'kitchen_qual', 'garage_sqft', 'garage_size', 'total_basement_sqft', 'sqft_1', 'basement_qual', 'year', 'garage_finish', 'garage_year', 'remod_year', 'baths', 'fireplace_qual', 'full_bath', 'mas_vnr_area', 'foundation', 'fireplaces', 'heating_qc', 'basement_exposure', 'basement_fin_sqft_1', 'gar_attached', 'sale_type', 'basement_fin_1', 'shape', 'fence', 'kitchen', 'conds', 'condition', 'contour_hill', 'contour_bank', 'contour_level', 'alley_gravel', 'alley_pave', 'street_material', 'lot_sqft', 'electrical', 'heating', 'basement_fin_2', 'sold_year_mo', 'basement_fin_sqft_2']"
$ws.Cells.Item(11, 5).Value2 = "full train"

# --- Row 12: ridge11_submit.csv ---
$ws.Cells.Item(12, 1).Value2 = "ridge11_submit.csv"
$ws.Cells.Item(12, 2).Value2 = "ridge"
$ws.Cells.Item(12, 3).Value2 = 900
$ws.Cells.Item(12, 4).Value2 = "['n_grnhill', 'n_greens', 'n_blueste', 'n_npkvill', 'n_veenker', 'n_brdale', 'n_blmngtn', 'n_meadowv', 'n_clearcr', 'n_swisu', 'n_stonebr', 'n_timber', 'n_noridge', 'n_idotrr', 'n_crawfor', 'n_brkside', 'n_mitchel', 'n_sawyerw', 'n_nwames', 'n_sawyer', 'n_gilbert', 'n_nridght', 'n_somerst', 'n_edwards', 'n_oldtown', 'n_collgcr', 'n_names', 'type_twn', 'type_sf', 'type_twn_end', 'gar_attached', 'gar_detached', 'gar_builtin', 'gar_basement', 'gar_2types', 'gar_carport', 'quality', 'gr_living_sqft', 'kitchen_qual', 'garage_sqft', 'garage_size', 'total_basement_sqft', 'sqft_1', 'basement_qual', 'year', 'garage_finish', 'garage_year', 'remod_year', 'baths', 'fireplace_qual', 'full_bath', 'mas_vnr_area', 'foundation', 'fireplaces', 'heating_qc', 'basement_exposure', 'basement_fin_sqft_1', 'gar_attached', 'sale_type', 'basement_fin_1', 'shape', 'fence', 'kitchen', 'conds', 'condition', 'contour_hill', 'contour_bank', 'contour_level', 'alley_gravel', 'alley_pave', 'street_material', 'lot_sqft', 'electrical', 'heating', 'basement_fin_2', 'sold_year_mo', 'basement_fin_sqft_2']"
$ws.Cells.Item(12, 5).Value2 = "full train"

# --- Row 13: ridge12_submit.csv (features cell wraps + adds low_qual_sqft) ---
# Features text entered before the filename, matching the authored edit order.
$ws.Cells.Item(13, 4).Value2 = "['n_grnhill', 'n_greens', 'n_blueste', 'n_npkvill', 'n_veenker', 'n_brdale', 'n_blmngtn', 'n_meadowv', 'n_clearcr', 'n_swisu', 'n_stonebr', 'n_timber', `n'n_noridge', 'n_idotrr', 'n_crawfor', 'n_brkside', 'n_mitchel', 'n_sawyerw', 'n_nwames', 'n_sawyer', 'n_gilbert', 'n_nridght', 'n_somerst', 'n_edwards', 'n_oldtown', `n'n_collgcr', 'n_names', 'type_twn', 'type_sf', 'type_twn_end', 'gar_attached', 'gar_detached', 'gar_builtin', 'gar_basement', 'gar_2types', 'gar_carport', 'quality', `n'gr_living_sqft', 'kitchen_qual', 'garage_sqft', 'garage_size', 'total_basement_sqft', 'sqft_1', 'basement_qual', 'year', 'garage_finish', 'garage_year', 'remod_year', `n'baths', 'fireplace_qual', 'full_bath', 'mas_vnr_area', 'foundation', 'fireplaces', 'heating_qc', 'basement_exposure', 'basement_fin_sqft_1', 'gar_attached', 'sale_type', `n'basement_fin_1', 'shape', 'fence', 'kitchen', 'conds', 'condition', 'contour_hill', 'contour_bank', 'contour_level', 'alley_gravel', 'alley_pave', 'street_material', `n'low_qual_sqft', 'lot_sqft', 'electrical', 'heating', 'basement_fin_2', 'sold_year_mo', 'basement_fin_sqft_2']"
$ws.Cells.Item(13, 4).WrapText = $true
$ws.Cells.Item(13, 1).Value2 = "ridge12_submit.csv"
$ws.Cells.Item(13, 2).Value2 = "ridge"
$ws.Cells.Item(13, 3).Value2 = 850
$ws.Cells.Item(13, 5).Value2 = "full train"
$ws.Rows.Item(13).RowHeight = 16

# --- Row 14: ridge13_submit.csv (features cell wraps, drops sold_year_mo/basement_fin_sqft_2) ---
$ws.Cells.Item(14, 4).Value2 = "['n_grnhill', 'n_greens', 'n_blueste', 'n_npkvill', 'n_veenker', 'n_brdale', 'n_blmngtn', 'n_meadowv', 'n_clearcr', 'n_swisu', 'n_stonebr', 'n_timber', `n'n_noridge', 'n_idotrr', 'n_crawfor', 'n_brkside', 'n_mitchel', 'n_sawyerw', 'n_nwames', 'n_sawyer', 'n_gilbert', 'n_nridght', 'n_somerst', 'n_edwards', 'n_oldtown', `n'n_collgcr', 'n_names', 'type_twn', 'type_sf', 'type_twn_end', 'gar_attached', 'gar_detached', 'gar_builtin', 'gar_basement', 'gar_2types', 'gar_carport', 'quality', `n'gr_living_sqft', 'kitchen_qual', 'garage_sqft', 'garage_size', 'total_basement_sqft', 'sqft_1', 'basement_qual', 'year', 'garage_finish', 'garage_year', 'remod_year', `n'baths', 'fireplace_qual', 'full_bath', 'mas_vnr_area', 'foundation', 'fireplaces', 'heating_qc', 'basement_exposure', 'basement_fin_sqft_1', 'gar_attached', 'sale_type', `n'basement_fin_1', 'shape', 'fence', 'kitchen', 'conds', 'condition', 'contour_hill', 'contour_bank', 'contour_level', 'alley_gravel', 'alley_pave', 'street_material', `n'low_qual_sqft', 'lot_sqft', 'electrical', 'heating', 'basement_fin_2']"
$ws.Cells.Item(14, 4).WrapText = $true
$ws.Cells.Item(14, 1).Value2 = "ridge13_submit.csv"
$ws.Cells.Item(14, 2).Value2 = "ridge"
$ws.Cells.Item(14, 3).Value2 = 800
$ws.Cells.Item(14, 5).Value2 = "full train"
$ws.Rows.Item(14).RowHeight = 15

# --- Row 15: ridge14_submit.csv (features cell wraps, drops low_qual_sqft) ---
$ws.Cells.Item(15, 1).Value2 = "ridge14_submit.csv"
$ws.Cells.Item(15, 2).Value2 = "ridge"
$ws.Cells.Item(15, 3).Value2 = 850
$ws.Cells.Item(15, 4).Value2 = "['n_grnhill', 'n_greens', 'n_blueste', 'n_npkvill', 'n_veenker', 'n_brdale', 'n_blmngtn', 'n_meadowv', 'n_clearcr', 'n_swisu', 'n_stonebr', 'n_timber', `n'n_noridge', 'n_idotrr', 'n_crawfor', 'n_brkside', 'n_mitchel', 'n_sawyerw', 'n_nwames', 'n_sawyer', 'n_gilbert', 'n_nridght', 'n_somerst', 'n_edwards', 'n_oldtown', 'n_collgcr', `n'n_names', 'type_twn', 'type_sf', 'type_twn_end', 'gar_attached', 'gar_detached', 'gar_builtin', 'gar_basement', 'gar_2types', 'gar_carport', 'quality', 'gr_living_sqft', `n'kitchen_qual', 'garage_sqft', 'garage_size', 'total_basement_sqft', 'sqft_1', 'basement_qual', 'year', 'garage_finish', 'garage_year', 'remod_year', 'baths', 'fireplace_qual', `n'full_bath', 'mas_vnr_area', 'foundation', 'fireplaces', 'heating_qc', 'basement_exposure', 'basement_fin_sqft_1', 'gar_attached', 'sale_type', 'basement_fin_1', 'shape', 'fence', `n'kitchen', 'conds', 'condition', 'contour_hill', 'contour_bank', 'contour_level', 'alley_gravel', 'alley_pave', 'street_material', 'lot_sqft', 'electrical', 'heating', `n'basement_fin_2', 'sold_year_mo', 'basement_fin_sqft_2']"
$ws.Cells.Item(15, 4).WrapText = $true
$ws.Cells.Item(15, 5).Value2 = "full train"
$ws.Rows.Item(15).RowHeight = 15

# Match the author's final selection/cursor position (last-edited cell).
$ws.Range("D15").Select()
